$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.058.42'
$ws.Range('E2').Value = '  -5.53%  '
$ws.Range('D3').Value = '2.225.70'
$ws.Range('E3').Value = '  -6.57%  '
$ws.Range('E4').Value = '  +0.22%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '323.09'
$ws.Range('E5').Value = '  -2.88%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '98.84'
$ws.Range('E6').Value = '  -9.69%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.581'
$ws.Range('E7').Value = '  -8.98%  '
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.563'
$ws.Range('E9').Value = '  -9.03%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '36.83'
$ws.Range('E10').Value = '  -10.55%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '54.29'
$ws.Range('E11').Value = '  -3.28%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0830'
$ws.Range('E12').Value = '  -10.10%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '7.64'
$ws.Range('E13').Value = '  -10.71%  '
$ws.Range('E14').Value = '  -2.13%  '
$ws.Range('B15').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C15').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D15').Value = '2.565.08'
$ws.Range('E15').Value = '  -6.29%  '
$ws.Range('B16').Value = 'Polygon'
$ws.Range('C16').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.862'
$ws.Range('E16').Value = '  -12.54%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '14.38'
$ws.Range('E17').Value = '  -7.60%  '
$ws.Range('D18').Value = '2.231.64'
$ws.Range('E18').Value = '  -6.09%  '
$ws.Range('D19').Value = '42.984.02'
$ws.Range('E19').Value = '  -5.52%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.01'
$ws.Range('E20').Value = '  -8.42%  '
$ws.Range('D21').Value = '0.0₃0967'
$ws.Range('E21').Value = '  -9.55%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.54'
$ws.Range('E22').Value = '  -11.03%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '3.22'
$ws.Range('E23').Value = '  -12.62%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '65.17'
$ws.Range('E24').Value = '  -11.33%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '237.26'
$ws.Range('E25').Value = '  -10.54%  '
$ws.Range('E26').Value = '  -6.80%  '
$ws.Range('E27').Value = '  -0.05%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.24'
$ws.Range('E29').Value = '  -2.93%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '9.99'
$ws.Range('E30').Value = '  -11.74%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.36'
$ws.Range('E31').Value = '  -15.31%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '36.28'
$ws.Range('E32').Value = '  -3.18%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '20.36'
$ws.Range('E33').Value = '  -9.71%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0868'
$ws.Range('E34').Value = '  -9.56%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '154.37'
$ws.Range('E35').Value = '  -8.82%  '
$ws.Range('E36').Value = '  -7.01%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.29'
$ws.Range('E37').Value = '  -1.03%  '
$ws.Range('E38').Value = '  -7.96%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.90'
$ws.Range('E39').Value = '  -4.50%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '4.42'
$ws.Range('E40').Value = '  -7.36%  '
$ws.Range('E41').Value = '  -11.74%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.71'
$ws.Range('E42').Value = '  -8.11%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0322'
$ws.Range('E43').Value = '  -9.76%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '13.89'
$ws.Range('E44').Value = '  +5.84%  '
$ws.Range('E45').Value = '  +0.08%  '
$ws.Range('D46').Value = '1.724.02'
$ws.Range('E46').Value = '  -8.28%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '84.99'
$ws.Range('E47').Value = '  -13.60%  '
$ws.Range('E48').Value = '  -12.06%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.95'
$ws.Range('E49').Value = '  -4.92%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '5.28'
$ws.Range('E50').Value = '  -13.45%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '74.85'
$ws.Range('E51').Value = '  -12.79%  '
